$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the placeholder "? # a".."? # e" labels (rows 36-40) to "? 1".."? 5"
$ws.Range("B36").Value = "? 1"
$ws.Range("C36").Value = "? 1"

$ws.Range("B37").Value = "? 1"
$ws.Range("C37").Value = "? 2"

$ws.Range("B38").Value = "? 1"
$ws.Range("C38").Value = "? 3"

$ws.Range("B39").Value = "? 1"
$ws.Range("C39").Value = "? 4"

$ws.Range("B40").Value = "? 1"
$ws.Range("C40").Value = "? 5"

# Reset the saved view: scroll back to the top and select A2 (instead of
# being scrolled down to A82 with B98 selected)
$ws.Range("A2").Select()
